$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MKTAY")

$updates = @{
    8 = @{ "D" = 4314800; "E" = 3751600; "F" = 3829600; "G" = 3749100; "H" = 3464200; "I" = 2799100; "J" = 2673200 }
    9 = @{ "D" = 2731600; "E" = 2418300; "F" = 2441900; "G" = 2328500; "H" = 2206200; "I" = 1761500; "J" = 1632100 }
    10 = @{ "D" = 1583100; "E" = 1333300; "F" = 1387700; "G" = 1420500; "H" = 1258000; "I" = 1037500; "J" = 1041100 }
    12 = @{ "D" = 98500; "E" = 91700; "F" = 86700; "G" = 82400; "H" = 78800; "I" = 75900 }
    17 = @{ "D" = 3593700; "E" = 3186000; "F" = 3244900; "G" = 3099000; "H" = 2967800; "I" = 2388900; "J" = 2234600 }
    18 = @{ "D" = 721000; "E" = 565600; "F" = 584700; "G" = 650000; "H" = 496400; "I" = 410100; "J" = 438600 }
    20 = @{ "E" = 20600; "F" = -27700; "G" = -30700; "H" = 20400; "I" = 4600; "J" = -11900 }
    21 = @{ "D" = 810400; "E" = 667100; "F" = 640900; "G" = 697300; "H" = 594900; "I" = 482900; "J" = 492200 }
    22 = @{ "J" = 2200 }
    23 = @{ "D" = 720300; "E" = 585200; "F" = 555900; "G" = 618300; "H" = 515000; "I" = 413000; "J" = 424500 }
    24 = @{ "D" = 220600; "E" = 177300; "F" = 176500; "G" = 205300; "H" = 164800; "I" = 131100; "J" = 128100 }
    26 = @{ "D" = 499600; "E" = 408000; "F" = 379400; "G" = 413000; "H" = 350200; "I" = 281900; "J" = 296400 }
    27 = @{ "D" = 495000; "E" = 404800; "F" = 376200; "G" = 409600; "H" = 347600; "I" = 280900; "J" = 293800 }
    32 = @{ "E" = -20600; "F" = 27700; "G" = 30700; "H" = -20400; "I" = -4600; "J" = 11900 }
    33 = @{ "D" = 495000; "E" = 404800; "F" = 376200; "G" = 409600; "H" = 347600; "I" = 280900; "J" = 293800 }
    35 = @{ "D" = 495000; "E" = 404800; "F" = 376200; "G" = 409600; "H" = 347600; "I" = 280900; "J" = 293800 }
    41 = @{ "D" = 1720600; "E" = 1595700; "F" = 1043800; "G" = 992700; "H" = 880500; "I" = 682900; "J" = 405100 }
    42 = @{ "D" = 567900; "E" = 211900; "F" = 436300; "G" = 506900; "H" = 371100; "I" = 344100; "J" = 349200 }
    43 = @{ "D" = 1429100; "E" = 607000; "F" = 580500; "G" = 587200; "H" = 583800; "I" = 488900; "J" = 447100 }
    44 = @{ "D" = 3549700; "E" = 1513300; "F" = 1616300; "G" = 1583700; "H" = 1411200; "I" = 1256100; "J" = 1171300 }
    45 = @{ "D" = 78300; "E" = 163400; "F" = 188400; "G" = 208600; "H" = 188900; "I" = 150400; "J" = 129200 }
    46 = @{ "D" = 4334300; "E" = 4091200; "F" = 3865300; "G" = 3879200; "H" = 3435500; "I" = 2922400; "J" = 2501900 }
    47 = @{ "D" = 713600; "E" = 307400; "F" = 197700; "G" = 283800; "H" = 274900; "I" = 166900; "J" = 173200 }
    48 = @{ "D" = 1814100; "E" = 854300; "F" = 839100; "G" = 866100; "H" = 828100; "I" = 780000; "J" = 702800 }
    49 = @{ "D" = 105400; "E" = 39400; "F" = 43600; "G" = 47800; "H" = 48900; "I" = 47600; "J" = 47300 }
    52 = @{ "D" = 203000; "E" = 106700; "F" = 98800; "G" = 124100; "H" = 105400; "I" = 69500; "J" = 39500 }
    54 = @{ "D" = 5919800; "E" = 5399100; "F" = 5044500; "G" = 5201000; "H" = 4692900; "I" = 3986400; "J" = 3464600 }
    57 = @{ "D" = 571100; "E" = 238200; "F" = 186400; "G" = 227100; "H" = 193500; "I" = 198100; "J" = 197300 }
    58 = @{ "D" = 65100; "E" = 59500; "F" = 19800; "G" = 42100; "H" = 37500; "I" = 15300; "J" = 21300 }
    59 = @{ "D" = 739700; "E" = 409900; "F" = 380400; "G" = 371400; "H" = 395400; "I" = 295800; "J" = 267100 }
    60 = @{ "D" = 733200; "E" = 707600; "F" = 586700; "G" = 640600; "H" = 626400; "I" = 509200; "J" = 485600 }
    61 = @{ "G" = 6900 }
    62 = @{ "D" = 166000; "E" = 118900; "F" = 87900; "G" = 127600; "H" = 93800; "I" = 75100; "J" = 52000 }
    66 = @{ "D" = 911200; "E" = 859500; "F" = 707600; "G" = 807300; "H" = 752000; "I" = 609600; "J" = 560500 }
    72 = @{ "D" = 8790100; "E" = 4178400; "F" = 3897500; "G" = 3666100; "H" = 3368200; "I" = 3108900; "J" = 2916400 }
    76 = @{ "D" = 5008600; "E" = 4539600; "F" = 4337000; "G" = 4393600; "H" = 3940800; "I" = 3376800; "J" = 2904100 }
    81 = @{ "D" = 495000; "E" = 404800; "F" = 376200; "G" = 409600; "H" = 347600; "I" = 280900; "J" = 293800 }
    83 = @{ "D" = 89600; "E" = 80800; "F" = 83900; "G" = 77900; "H" = 77900; "I" = 68200; "J" = 65400 }
    89 = @{ "D" = 301000; "E" = 572700; "F" = 309100; "G" = 324500; "H" = 376800; "I" = 346800; "J" = 77900 }
    91 = @{ "D" = -127900; "E" = -119500; "F" = -106400; "G" = -109600; "H" = -103200; "I" = -103800; "J" = -121900 }
    94 = @{ "D" = -133900; "E" = -45300; "F" = -59400; "G" = -181700; "H" = -181600; "I" = -139300; "J" = -40700 }
    96 = @{ "D" = -125200; "E" = -123900; "F" = -144800; "G" = -111700; "H" = -88300; "I" = -88400; "J" = -82200 }
    100 = @{ "D" = -160400; "E" = -85800; "F" = -169200; "G" = -108600; "H" = -66600; "I" = -96300; "J" = -114900 }
    101 = @{ "D" = 39800; "E" = -59400; "F" = -31700; "G" = 81500; "H" = 47100; "I" = 46700; "J" = 14100 }
    102 = @{ "D" = 46500; "E" = 382100; "F" = 48700; "G" = 115700; "H" = 175800; "I" = 157900; "J" = -63500 }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
